# Add changes for attachments: new "mail_attachment_title" header column (G)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cell in G1 -> introduces the new shared string "mail_attachment_title"
$ws.Range("G1").Value = "mail_attachment_title"

# Give the new column G a custom width (matches the author's saved layout)
$ws.Columns.Item(7).ColumnWidth = 25

# Move/save the active selection to G4, like in the committed workbook
$ws.Range("G4").Select()
